$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q2: new order number, and the cell loses its centered style (now default/unstyled)
$ws.Range("Q2").ClearFormats()
$ws.Range("Q2").Value = 4503342107

# Q3: new order number, style (centered) is kept
$ws.Range("Q3").Value = 4503342108

# Update the active selection to P4 (was P5)
$ws.Range("P4").Select()
